# Update cryptocurrency price/volume figures per the upstream GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.210.46"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "'1.904.00"
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("E5").Value = "  +0.25%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.21%  "
$ws.Range("E7").Value = "  -0.81%  "
$ws.Range("D8").Value = "'0.3818"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.98%  "
$ws.Range("D9").Value = "'0.07312"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.83%  "
$ws.Range("D10").Value = "'21.72"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.57%  "
$ws.Range("D11").Value = "'0.9057"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.56%  "
$ws.Range("D12").Value = "'0.08069"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.86%  "
$ws.Range("D13").Value = "'95.68"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.71%  "
$ws.Range("D14").Value = "'5.356"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.55%  "
$ws.Range("D15").Value = "'1.824.50"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.23%  "
$ws.Range("D16").Value = "'1.003"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.29%  "
$ws.Range("D17").Value = "'0.000008683"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.73%  "
$ws.Range("E18").Value = "  +1.02%  "
$ws.Range("E19").Value = "  +0.15%  "
$ws.Range("D20").Value = "'27.243.06"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'5.118"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.01%  "
$ws.Range("E22").Value = "  +1.96%  "
$ws.Range("D23").Value = "'6.473"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.44%  "
$ws.Range("D24").Value = "'2.343"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.65%  "
$ws.Range("D25").Value = "'149.39"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.36%  "
$ws.Range("E26").Value = "  +0.35%  "
$ws.Range("E27").Value = "  -0.49%  "
$ws.Range("D28").Value = "'116.14"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.00%  "
$ws.Range("D29").Value = "'4.837"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.18%  "
$ws.Range("D30").Value = "'4.899"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Value = "'0.09238"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.84%  "
$ws.Range("D32").Value = "'0.05077"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.14%  "
$ws.Range("D33").Value = "'0.7999"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.37%  "
$ws.Range("D34").Value = "'1.228"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.91%  "
$ws.Range("D35").Value = "'2.978"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.92%  "
$ws.Range("D36").Value = "'3.372"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.12%  "
$ws.Range("D37").Value = "'2.666"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.55%  "
$ws.Range("D38").Value = "'0.5732"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("D39").Value = "'0.01994"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.24%  "
$ws.Range("D40").Value = "'1.087"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.14%  "
$ws.Range("D41").Value = "'9.010"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.25%  "
$ws.Range("D42").Value = "'6.600"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.74%  "
$ws.Range("D43").Value = "'116.65"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.83%  "
$ws.Range("D44").Value = "'0.1518"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.16%  "
$ws.Range("D45").Value = "'0.4900"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.00%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'10.20"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.09%  "
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").Value = "'1.002"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.22%  "
$ws.Range("D48").Value = "'1.635"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.05%  "
$ws.Range("D49").Value = "'38.58"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.93%  "
$ws.Range("D50").Value = "'64.17"
$ws.Range("D50").Style = "Normal"
$ws.Range("E51").Value = "  +0.39%  "
